$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2 (with formatting) down to rows 3-11 so the new rows reuse the same cell styles
$srcRow = $ws.Range("A2:Q2")
$srcRow.Copy($ws.Range("A3:Q3"))
$srcRow.Copy($ws.Range("A4:Q4"))
$srcRow.Copy($ws.Range("A5:Q5"))
$srcRow.Copy($ws.Range("A6:Q6"))
$srcRow.Copy($ws.Range("A7:Q7"))
$srcRow.Copy($ws.Range("A8:Q8"))
$srcRow.Copy($ws.Range("A9:Q9"))
$srcRow.Copy($ws.Range("A10:Q10"))
$srcRow.Copy($ws.Range("A11:Q11"))

# Update dimension-driving data for row 2 and the newly added rows 3-11 (monthly -> weekly cadence)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 44381
$ws.Range("C2").Value = 695.9751779659155
$ws.Range("D2").Value = 421.9468285179482
$ws.Range("E2").Value = 647.7038520383048
$ws.Range("F2").Value = 695.9751778747167
$ws.Range("G2").Value = 695.9751780347525
$ws.Range("H2").Value = -155.7333778341101
$ws.Range("I2").Value = -155.7333778341101
$ws.Range("J2").Value = -155.7333778341101
$ws.Range("K2").Value = -155.7333778341101
$ws.Range("L2").Value = -155.7333778341101
$ws.Range("M2").Value = -155.7333778341101
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 540.2418001318053

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 44388
$ws.Range("C3").Value = 698.0403831983502
$ws.Range("D3").Value = 513.0300394466716
$ws.Range("E3").Value = 748.4921605809195
$ws.Range("F3").Value = 698.0403829056385
$ws.Range("G3").Value = 698.0403834659427
$ws.Range("H3").Value = -65.67529490363245
$ws.Range("I3").Value = -65.67529490363245
$ws.Range("J3").Value = -65.67529490363245
$ws.Range("K3").Value = -65.67529490363245
$ws.Range("L3").Value = -65.67529490363245
$ws.Range("M3").Value = -65.67529490363245
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 632.3650882947178

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 44395
$ws.Range("C4").Value = 700.105588430785
$ws.Range("D4").Value = 683.3157487796057
$ws.Range("E4").Value = 892.0426671812799
$ws.Range("F4").Value = 700.1055878724994
$ws.Range("G4").Value = 700.1055889578307
$ws.Range("H4").Value = 82.94745979380299
$ws.Range("I4").Value = 82.94745979380299
$ws.Range("J4").Value = 82.94745979380299
$ws.Range("K4").Value = 82.94745979380299
$ws.Range("L4").Value = 82.94745979380299
$ws.Range("M4").Value = 82.94745979380299
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 783.053048224588

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 44402
$ws.Range("C5").Value = 702.1707936632197
$ws.Range("D5").Value = 827.9521295639863
$ws.Range("E5").Value = 1061.146012018548
$ws.Range("F5").Value = 702.1707927895143
$ws.Range("G5").Value = 702.1707945046884
$ws.Range("H5").Value = 240.3593818018295
$ws.Range("I5").Value = 240.3593818018295
$ws.Range("J5").Value = 240.3593818018295
$ws.Range("K5").Value = 240.3593818018295
$ws.Range("L5").Value = 240.3593818018295
$ws.Range("M5").Value = 240.3593818018295
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 942.5301754650492

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 44409
$ws.Range("C6").Value = 704.2359988956543
$ws.Range("D6").Value = 933.8293929072228
$ws.Range("E6").Value = 1163.36670002597
$ws.Range("F6").Value = 704.2359976535694
$ws.Range("G6").Value = 704.2360001099456
$ws.Range("H6").Value = 347.0269879323832
$ws.Range("I6").Value = 347.0269879323832
$ws.Range("J6").Value = 347.0269879323832
$ws.Range("K6").Value = 347.0269879323832
$ws.Range("L6").Value = 347.0269879323832
$ws.Range("M6").Value = 347.0269879323832
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1051.262986828038

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 44416
$ws.Range("C7").Value = 706.3012041280891
$ws.Range("D7").Value = 965.479476789194
$ws.Range("E7").Value = 1196.441127546732
$ws.Range("F7").Value = 706.3012024482583
$ws.Range("G7").Value = 706.3012057388836
$ws.Range("H7").Value = 377.0436211460323
$ws.Range("I7").Value = 377.0436211460323
$ws.Range("J7").Value = 377.0436211460323
$ws.Range("K7").Value = 377.0436211460323
$ws.Range("L7").Value = 377.0436211460323
$ws.Range("M7").Value = 377.0436211460323
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1083.344825274121

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 44423
$ws.Range("C8").Value = 708.3664093605238
$ws.Range("D8").Value = 933.0787167339626
$ws.Range("E8").Value = 1159.203736319341
$ws.Range("F8").Value = 708.3664071476358
$ws.Range("G8").Value = 708.3664113754444
$ws.Range("H8").Value = 342.9768708091839
$ws.Range("I8").Value = 342.9768708091839
$ws.Range("J8").Value = 342.9768708091839
$ws.Range("K8").Value = 342.9768708091839
$ws.Range("L8").Value = 342.9768708091839
$ws.Range("M8").Value = 342.9768708091839
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1051.343280169708

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 44430
$ws.Range("C9").Value = 710.4316145929585
$ws.Range("D9").Value = 863.3868951005149
$ws.Range("E9").Value = 1086.660223786152
$ws.Range("F9").Value = 710.4316118603862
$ws.Range("G9").Value = 710.4316170765194
$ws.Range("H9").Value = 264.8100394333664
$ws.Range("I9").Value = 264.8100394333664
$ws.Range("J9").Value = 264.8100394333664
$ws.Range("K9").Value = 264.8100394333664
$ws.Range("L9").Value = 264.8100394333664
$ws.Range("M9").Value = 264.8100394333664
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 975.2416540263248

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 44437
$ws.Range("C10").Value = 712.4968198253932
$ws.Range("D10").Value = 744.7456466478602
$ws.Range("E10").Value = 968.6029353028998
$ws.Range("F10").Value = 712.496816562929
$ws.Range("G10").Value = 712.4968227957229
$ws.Range("H10").Value = 146.5658959607051
$ws.Range("I10").Value = 146.5658959607051
$ws.Range("J10").Value = 146.5658959607051
$ws.Range("K10").Value = 146.5658959607051
$ws.Range("L10").Value = 146.5658959607051
$ws.Range("M10").Value = 146.5658959607051
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 859.0627157860984

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 44444
$ws.Range("C11").Value = 714.562025057828
$ws.Range("D11").Value = 592.1117150714443
$ws.Range("E11").Value = 815.6346574254367
$ws.Range("F11").Value = 714.5620211961696
$ws.Range("G11").Value = 714.5620285998358
$ws.Range("H11").Value = -10.45737897078865
$ws.Range("I11").Value = -10.45737897078865
$ws.Range("J11").Value = -10.45737897078865
$ws.Range("K11").Value = -10.45737897078865
$ws.Range("L11").Value = -10.45737897078865
$ws.Range("M11").Value = -10.45737897078865
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 704.1046460870393

